$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header description cell (new run date / event count)
$ws.Range("A1").Value = "Description unknown, completed 06/21/2023 08:45:38 EDT, by WPJTOWN1.The search returned: 9 events."

# --- Row 4 (now collapses to an "Not authorized" style row, like row 3) ---
$ws.Range("A4").Value = "CGAX"
$ws.Range("B4").Value = 10266
$ws.Range("C4").Value = "Not authorized to view shipment"
$ws.Range("D4").Value = $null
$ws.Range("E4").Value = $null
$ws.Range("F4").Value = $null
$ws.Range("G4").Value = $null
$ws.Range("H4").Value = $null
$ws.Range("I4").Value = $null
$ws.Range("J4").Value = $null
$ws.Range("K4").Value = $null
$ws.Range("L4").Value = "Not authorized to view shipment"
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = $null
$ws.Range("O4").Value = "CGAX10266"

# --- Row 5 ---
$ws.Range("A5").Value = "CGAX"
$ws.Range("B5").Value = 10087
$ws.Range("C5").Value = "DODGE CITY"
$ws.Range("D5").Value = "KS"
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 20
$ws.Range("G5").Value = 2209
$ws.Range("H5").Value = "Arrive In-Transit"
$ws.Range("I5").Value = "LKAN01"
$ws.Range("J5").Value = "LOVELAND"
$ws.Range("K5").Value = "CO"
$ws.Range("L5").Value = 260040
$ws.Range("M5").Value = 63900
$ws.Range("N5").Value = 196140
$ws.Range("O5").Value = "CGAX10087"

# --- Row 6 ---
$ws.Range("A6").Value = "CRDX"
$ws.Range("B6").Value = 15803
$ws.Range("C6").Value = "JOHNSTOWN"
$ws.Range("D6").Value = "CO"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 15
$ws.Range("G6").Value = 1435
$ws.Range("H6").Value = "Placed Actual"
$ws.Range("I6").Value = $null
$ws.Range("J6").Value = "LOVELAND"
$ws.Range("K6").Value = "CO"
$ws.Range("L6").Value = 284700
$ws.Range("M6").Value = 66900
$ws.Range("N6").Value = 217800
$ws.Range("O6").Value = "CRDX15803"

# --- Row 7 ---
$ws.Range("A7").Value = "HRTX"
$ws.Range("B7").Value = 541048
$ws.Range("C7").Value = "JOHNSTOWN"
$ws.Range("D7").Value = "CO"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 15
$ws.Range("G7").Value = 1435
$ws.Range("H7").Value = "Placed Actual"
$ws.Range("I7").Value = $null
$ws.Range("J7").Value = "LOVELAND"
$ws.Range("K7").Value = "CO"
$ws.Range("L7").Value = 202800
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 202800
$ws.Range("O7").Value = "HRTX541048"

# --- Row 8 ---
$ws.Range("A8").Value = "CRDX"
$ws.Range("B8").Value = 15008
$ws.Range("C8").Value = "JOHNSTOWN"
$ws.Range("D8").Value = "CO"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 17
$ws.Range("G8").Value = 1431
$ws.Range("H8").Value = "Placed Actual"
$ws.Range("I8").Value = $null
$ws.Range("J8").Value = "LOVELAND"
$ws.Range("K8").Value = "CO"
$ws.Range("L8").Value = 286650
$ws.Range("M8").Value = 68700
$ws.Range("N8").Value = 217950
$ws.Range("O8").Value = "CRDX15008"

# --- Row 9 ---
$ws.Range("A9").Value = "HRTX"
$ws.Range("B9").Value = 541043
$ws.Range("C9").Value = "KANSAS CITY"
$ws.Range("D9").Value = "KS"
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 19
$ws.Range("G9").Value = 2345
$ws.Range("H9").Value = "Arrive In-Transit"
$ws.Range("I9").Value = "HLINKC"
$ws.Range("J9").Value = "LOVELAND"
$ws.Range("K9").Value = "CO"
$ws.Range("L9").Value = 258850
$ws.Range("M9").Value = 64200
$ws.Range("N9").Value = 194650
$ws.Range("O9").Value = "HRTX541043"

# --- Row 10 (new) ---
$ws.Range("A10").Value = "HRTX"
$ws.Range("B10").Value = 541059
$ws.Range("C10").Value = "LOVELAND"
$ws.Range("D10").Value = "CO"
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 19
$ws.Range("G10").Value = 1215
$ws.Range("H10").Value = "Junction Received"
$ws.Range("I10").Value = "BNSF"
$ws.Range("J10").Value = "LOVELAND"
$ws.Range("K10").Value = "CO"
$ws.Range("L10").Value = 261250
$ws.Range("M10").Value = 64200
$ws.Range("N10").Value = 197050
$ws.Range("O10").Value = "HRTX541059"

# --- Row 11 (new) ---
$ws.Range("A11").Value = "BNSF"
$ws.Range("B11").Value = 468933
$ws.Range("C11").Value = "LOVELAND"
$ws.Range("D11").Value = "CO"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 19
$ws.Range("G11").Value = 1215
$ws.Range("H11").Value = "Junction Received"
$ws.Range("I11").Value = "BNSF"
$ws.Range("J11").Value = "LOVELAND"
$ws.Range("K11").Value = "CO"
$ws.Range("L11").Value = 234960
$ws.Range("M11").Value = 63600
$ws.Range("N11").Value = 171360
$ws.Range("O11").Value = "BNSF468933"

# Update selection to match new extended range (O3:O11)
$ws.Range("O3:O11").Select()
